$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text: "No." -> "No" and "Lokasi" -> "Lokasi Barang"
$ws.Range("A1").Value = "No"
$ws.Range("G1").Value = "Lokasi Barang"

# Match the selection shown in the post-edit file
$ws.Range("C6").Select()
